# The document ends with a "Requisitos" section followed by some
# footer-like content that was scraped from the site chrome ("Ver no
# Jupiter Salvar em pdf Salvar em docx" and the "© 2020 ..." copyright
# line), preceded by a blank paragraph. This footer content (and its
# leading blank paragraph) needs to be removed, while keeping the blank
# paragraph that sits right before the trailing page-break paragraph.
#
# Locate the three paragraphs to remove by their text and delete the
# Range spanning from the start of the first one through the end
# (including the paragraph mark) of the last one.

$d = $word.ActiveDocument

$startPara = $null
$endPara = $null

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text

    if ($text -like "*Ver no Jupiter*") {
        # The blank paragraph immediately preceding this one is also
        # part of the block to remove.
        $startPara = $d.Paragraphs.Item($i - 1)
    }
    if ($text -like "*Powered by Jekyll*") {
        $endPara = $para
    }
}

if ($startPara -ne $null -and $endPara -ne $null) {
    $rng = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $rng.Delete()
}
